$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of G2 and G3 (Job values for the two rows)
$ws.Range("G2").Value = "L3 (Expert)"
$ws.Range("G3").Value = "L4 (Professional)"

# Update L3 (Date of Birth) to a new date value (serial 31982 -> 1987-07-24)
$ws.Range("L3").Value = (Get-Date -Year 1987 -Month 7 -Day 24 -Hour 0 -Minute 0 -Second 0)

# Update the selection shown in the sheet view
$ws.Range("G5:G8").Select()
